$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two card-number cells were re-entered as text (last digit changed) --
# force text storage so Excel keeps the leading/looks-like-a-number string
# instead of re-parsing it back into a numeric value.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "4595980021202763"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "5892290000424181"

# Suppress the "number stored as text" warning triangle on those cells.
$ws.Range("B4:B5").Errors.Item(3).Ignore = $true

# Selection moved to B10 before the file was saved.
$ws.Range("B10").Select() | Out-Null
